# Facilitator guidelines - Playful mathematicians.docx
# English -> Swahili (Kenya) translation edits.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute(
        $old,      # FindText
        $true,     # MatchCase
        $true,     # MatchWholeWord
        $false,    # MatchWildcards
        $false,    # MatchSoundsLike
        $false,    # MatchAllWordForms
        $true,     # Forward
        1,         # Wrap (wdFindContinue)
        $false,    # Format
        $new,      # ReplaceWith
        2          # Replace (wdReplaceAll)
    ) | Out-Null
}

Replace-Exact "Video Title" "Kichwa cha Video"
Replace-Exact "Topic" "Mada"
Replace-Exact "Aim(s)" "Malengo"
Replace-Exact "Length" "Urefu"
Replace-Exact "Camp Location" "Mahali pa Kambi"
Replace-Exact "Facilitators" "Wawezeshaji"
Replace-Exact "N. of students" "N. ya wanafunzi"
Replace-Exact "Date" "Tarehe"
Replace-Exact "Resources" "Rasilimali"
Replace-Exact "needed" "inahitajika"
Replace-Exact "Preparations" "Maandalizi"
Replace-Exact "Video time" "Muda wa video"
Replace-Exact "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Exact "What learners do" "Wanachofanya wanafunzi"
Replace-Exact "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-Exact "Video Introduction" "Utangulizi wa Video"
Replace-Exact "Riddle" "Kitendawili"
Replace-Exact "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
Replace-Exact "Solution" "Suluhisho"

# Document default language: Swahili (Tanzania) -> Swahili (Kenya).
# The document-wide default lives in styles.xml docDefaults, which isn't
# directly addressable through the exposed Word object model; setting the
# language on the (default) "Normal" style's font achieves the same
# effective default-language result for all content in this document
# (every paragraph in this document uses the Normal style).
$d.Styles("Normal").Font.LanguageID = "sw-KE"

"OK"
